$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. Columns B:E (values) shift to C:F,
# the new column B inherits column A's style and is currently empty,
# while column A (with the segment names) is left untouched.
$ws.Columns.Item(2).Insert()

# New header for the inserted column; copy the header style (bold,
# bordered, centered) from the neighbouring header cell since the
# inserted column picked up column A's (unstyled) formatting in row 1.
$ws.Range("B1").Value = "segments"
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# For each data row: move the segment name from column A into the new
# column B, and replace column A with a 0-based numeric index.
for ($row = 2; $row -le 20; $row++) {
    $name = $ws.Cells.Item($row, 1).Value()
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 1).Value = $row - 2
}

# The inserted column B picked up column A's bold/bordered/centred data
# style; the segment-name column should stay unstyled like the other
# value columns, so strip the inherited formatting from the data rows.
$ws.Range("B2:B20").ClearFormats()
